$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. Fix capitalization of the English footnote shared string in C7
#    ("*according ..." -> "*According ...")
# ---------------------------------------------------------------------------
$ws.Range("C7").Value = "*According to the Service for the Regulation and Supervision of the Communications Sector under the Ministry of Digital Development of the Kyrgyz Republic"

# ---------------------------------------------------------------------------
# 2. Re-map the unused custom number format (164) on J7 to the equivalent
#    built-in accounting format (41) by re-applying the identical format
#    code; the engine recognizes it matches a built-in format and reuses it.
# ---------------------------------------------------------------------------
$ws.Range("J7").NumberFormat = "_(* #,##0_);_(* (#,##0);_(* ""-""_);_(@_)"

# ---------------------------------------------------------------------------
# 3. Add the new 2023 column (O) mirroring column N's layout/styles.
# ---------------------------------------------------------------------------
$ws.Range("N2").Copy() | Out-Null
$ws.Range("O2").PasteSpecial(-4122) | Out-Null

$ws.Range("N3").Copy() | Out-Null
$ws.Range("O3").PasteSpecial(-4122) | Out-Null
$ws.Range("O3").Value = 2023

$ws.Range("N4").Copy() | Out-Null
$ws.Range("O4").PasteSpecial(-4122) | Out-Null
$ws.Range("O4").Value = 5571

$ws.Range("N5").Copy() | Out-Null
$ws.Range("O5").PasteSpecial(-4122) | Out-Null
$ws.Range("O5").Value = 74710

$ws.Range("N6").Copy() | Out-Null
$ws.Range("O6").PasteSpecial(-4122) | Out-Null
$ws.Range("O6").Value = 375715

# ---------------------------------------------------------------------------
# 4. Row height tweaks.
# ---------------------------------------------------------------------------
$ws.Rows.Item(1).RowHeight = 45
$ws.Rows.Item(7).RowHeight = 46.5

# ---------------------------------------------------------------------------
# 5. Shrink the font used for the footnote/legend row (row 7) from size 9/10
#    down to size 8, reusing the formatting already present on row 4 as a
#    base (so only the font size actually changes).
# ---------------------------------------------------------------------------
$ws.Range("A4").Copy() | Out-Null
$ws.Range("A7").PasteSpecial(-4122) | Out-Null
$ws.Range("A7").Font.Size = 8

$ws.Range("B4").Copy() | Out-Null
$ws.Range("B7:C7").PasteSpecial(-4122) | Out-Null
$ws.Range("B7:C7").Font.Size = 8

$excel.CutCopyMode = 0
